# Apply the "automatic update" of the Artfynd export data.
# Rows 7-10 got refreshed with new Id/Ost/Nord coordinates, a new
# Taxonsorteringsordning (column B) value, and the "Antal" (column I)
# values moved along with their originating records.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 7 ---
$ws.Range("A7").Value = 112092161
$ws.Range("B7").Value = 96720
$ws.Range("I7").Value = "'10"
$ws.Range("Q7").Value = 584330
$ws.Range("R7").Value = 7048274

# --- Row 8 ---
$ws.Range("A8").Value = 112092066
$ws.Range("B8").Value = 96720
$ws.Range("I8").Value = ""
$ws.Range("Q8").Value = 584346
$ws.Range("R8").Value = 7048207
$ws.Range("Z8").Value = "17:18"
$ws.Range("AB8").Value = "17:18"

# --- Row 9 ---
$ws.Range("A9").Value = 112092130
$ws.Range("B9").Value = 96720
$ws.Range("I9").Value = ""
$ws.Range("Q9").Value = 584352
$ws.Range("R9").Value = 7048232
$ws.Range("Z9").Value = "17:22"
$ws.Range("AB9").Value = "17:22"

# --- Row 10 ---
$ws.Range("A10").Value = 112092586
$ws.Range("B10").Value = 96720
$ws.Range("I10").Value = "'20"
$ws.Range("Q10").Value = 584401
$ws.Range("R10").Value = 7048357
$ws.Range("Z10").Value = "17:46"
$ws.Range("AB10").Value = "17:46"
